# TrackAttendance.xlsx - "updated by week 3 - class 6"
# Marks attendance (0 -> 1) for class 6 (column G, "aula6") and, where it
# had been missed previously, also fills in class 5 (column F, "aula5")
# and in a few rows class 4 (column E, "aula4") so the running attendance
# stays contiguous. Also highlights the still-blank roster row (44) in
# yellow and leaves the selection on the last-edited cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cellsToMark = @(
    "F2","G2","F4","G4","G5","E6","G6","F7","G7","F8","G8","F9","G9",
    "F10","G10","F11","G11","G12","F13","F14","G14","F15","G15","F16",
    "G16","F17","G17","F18","E19","F19","G19","F20","G20","F21","G21",
    "F22","F23","G23","F24","G24","G25","F26","G26","F27","G27","F28",
    "G28","F29","G29","F30","G30","F31","G31","F32","G32","F33","G33",
    "F34","G34","F35","G35","F36","G36","E37","F37","G37","F38","G38",
    "F39","G39","F40","G40","F41","F42","G42","F43","G43","F45","G45",
    "E46","F46","G46","F47","F48","G48","F49","G49"
)

foreach ($addr in $cellsToMark) {
    $ws.Range($addr).Value = 1
}

# Highlight the empty row (row 44) in yellow, same as Excel's standard
# "Fill Color -> Yellow" applied to B44:Y44.
$ws.Range("B44:Y44").Interior.Color = 65535

# Leave the selection where the editor last clicked.
[void]$ws.Range("AD22").Select()
